$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values changed (all originally stored as text/shared-strings,
# representing numeric-looking metric values). We force the "@" (text) number
# format before assigning so Excel keeps them as text rather than converting
# to numbers, then clear the format again so no stray style lingers on the
# cell (this also prevents the format from spreading to other cells).

$textChanges = @{
    "G2" = "69"
    "G3" = "69"
    "G4" = "69"
    "G5" = "69"
    "G6" = "69"
    "G7" = "69"
    "G8" = "69"
    "G9" = "13"
    "J4" = "11"
    "J5" = "11"
    "J6" = "11"
    "J7" = "11"
    "J8" = "12"
    "J9" = "13"
}

foreach ($addr in $textChanges.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $textChanges[$addr]
    $range.ClearFormats()
}

